{"js": "// Replace each old three-digit-by-one-digit multiplication expression\n// with its new counterpart, per the commit diff.\nconst replacements = [\n  [\"611\u00d78=4888\", \"769\u00d73=2307\"],\n  [\"208\u00d78=1664\", \"571\u00d76=3426\"],\n  [\"105\u00d75=525\", \"930\u00d79=8370\"],\n  [\"272\u00d73=816\", \"211\u00d79=1899\"],\n  [\"557\u00d77=3899\", \"410\u00d75=2050\"],\n  [\"717\u00d74=2868\", \"136\u00d73=408\"],\n  [\"973\u00d74=3892\", \"243\u00d74=972\"],\n  [\"851\u00d73=2553\", \"214\u00d78=1712\"],\n  [\"835\u00d74=3340\", \"292\u00d74=1168\"],\n  [\"134\u00d78=1072\", \"413\u00d79=3717\"],\n  [\"433\u00d74=1732\", \"171\u00d76=1026\"],\n  [\"904\u00d73=2712\", \"873\u00d78=6984\"],\n  [\"748\u00d78=5984\", \"572\u00d76=3432\"],\n  [\"369\u00d78=2952\", \"157\u00d78=1256\"],\n  [\"641\u00d77=4487\", \"117\u00d74=468\"],\n  [\"566\u00d75=2830\", \"608\u00d79=5472\"],\n  [\"102\u00d76=612\", \"617\u00d75=3085\"],\n  [\"178\u00d72=356\", \"572\u00d78=4576\"],\n  [\"105\u00d74=420\", \"919\u00d76=5514\"],\n  [\"949\u00d72=1898\", \"716\u00d75=3580\"],\n  [\"909\u00d78=7272\", \"661\u00d77=4627\"],\n  [\"989\u00d76=5934\", \"382\u00d79=3438\"],\n  [\"714\u00d75=3570\", \"544\u00d78=4352\"],\n  [\"572\u00d79=5148\", \"531\u00d72=1062\"],\n  [\"219\u00d77=1533\", \"926\u00d76=5556\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each old three-digit-by-one-digit multiplication expression\n# with its new counterpart, per the commit diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{old=\"611\u00d78=4888\"; new=\"769\u00d73=2307\"},\n    @{old=\"208\u00d78=1664\"; new=\"571\u00d76=3426\"},\n    @{old=\"105\u00d75=525\"; new=\"930\u00d79=8370\"},\n    @{old=\"272\u00d73=816\"; new=\"211\u00d79=1899\"},\n    @{old=\"557\u00d77=3899\"; new=\"410\u00d75=2050\"},\n    @{old=\"717\u00d74=2868\"; new=\"136\u00d73=408\"},\n    @{old=\"973\u00d74=3892\"; new=\"243\u00d74=972\"},\n    @{old=\"851\u00d73=2553\"; new=\"214\u00d78=1712\"},\n    @{old=\"835\u00d74=3340\"; new=\"292\u00d74=1168\"},\n    @{old=\"134\u00d78=1072\"; new=\"413\u00d79=3717\"},\n    @{old=\"433\u00d74=1732\"; new=\"171\u00d76=1026\"},\n    @{old=\"904\u00d73=2712\"; new=\"873\u00d78=6984\"},\n    @{old=\"748\u00d78=5984\"; new=\"572\u00d76=3432\"},\n    @{old=\"369\u00d78=2952\"; new=\"157\u00d78=1256\"},\n    @{old=\"641\u00d77=4487\"; new=\"117\u00d74=468\"},\n    @{old=\"566\u00d75=2830\"; new=\"608\u00d79=5472\"},\n    @{old=\"102\u00d76=612\"; new=\"617\u00d75=3085\"},\n    @{old=\"178\u00d72=356\"; new=\"572\u00d78=4576\"},\n    @{old=\"105\u00d74=420\"; new=\"919\u00d76=5514\"},\n    @{old=\"949\u00d72=1898\"; new=\"716\u00d75=3580\"},\n    @{old=\"909\u00d78=7272\"; new=\"661\u00d77=4627\"},\n    @{old=\"989\u00d76=5934\"; new=\"382\u00d79=3438\"},\n    @{old=\"714\u00d75=3570\"; new=\"544\u00d78=4352\"},\n    @{old=\"572\u00d79=5148\"; new=\"531\u00d72=1062\"},\n    @{old=\"219\u00d77=1533\"; new=\"926\u00d76=5556\"}\n)\n\nforeach ($p in $pairs) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $p.old\n    $find.Replacement.Text = $p.new\n    $found = $find.Execute(\n        $p.old,      # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $p.new,      # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n    if (-not $found) {\n        throw \"Text not found: $($p.old)\"\n    }\n}\n"}
